$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp (row 1): 08:04 -> 09:04
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 09:04"

# 2) Estados Unidos (row 4)
$ws.Range("B4").Value = 1347318
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 238080
$ws.Range("E4").Value = 1029198
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 80040

# 3) Alemania (row 10)
$ws.Range("D10").Value = 144400
$ws.Range("E10").Value = 19375

# 4) Georgia (row 119)
$ws.Range("B119").Value = 635
$ws.Range("C119").Value = 9
$ws.Range("D119").Value = 309
$ws.Range("E119").Value = 316

# 5) Taiwan (row 125)
$ws.Range("D125").Value = 366
$ws.Range("E125").Value = 68

# 6) Swap the Belice / Nueva Caledonia rows (192 <-> 193), full row contents
$ws.Range("A192").Value = "Belice"
$ws.Range("B192").Value = 18
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 16
$ws.Range("E192").Value = 0
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 2

$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("B193").Value = 18
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 18
$ws.Range("E193").Value = 0
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 0
